$wb = $excel.ActiveWorkbook

# --- Sheet: "Trends Status" ---
$ws1 = $wb.Worksheets.Item("Trends Status")

$ws1.Range("B2").Value = 38
$ws1.Range("C2").Value = 34
$ws1.Range("D2").Value = 16.7
$ws1.Range("E2").Value = 11.8

$ws1.Range("B3").Value = 59
$ws1.Range("C3").Value = 75
$ws1.Range("D3").Value = 25.9
$ws1.Range("E3").Value = 26

$ws1.Range("B4").Value = 74
$ws1.Range("C4").Value = 139
$ws1.Range("D4").Value = 32.5
$ws1.Range("E4").Value = 48.3

$ws1.Range("B5").Value = 26
$ws1.Range("C5").Value = 16
$ws1.Range("D5").Value = 11.4
$ws1.Range("E5").Value = 5.6

$ws1.Range("B6").Value = 31
$ws1.Range("C6").Value = 24
$ws1.Range("D6").Value = 13.6
$ws1.Range("E6").Value = 8.300000000000001

$ws1.Range("B7").Value = 298
$ws1.Range("C7").Value = 363

# --- Sheet: "Interannual update - Reason _1" ---
$ws12 = $wb.Worksheets.Item("Interannual update - Reason _1")

$ws12.Range("B2").Value = 13
$ws12.Range("C2").Value = 16.2

$ws12.Range("B3").Value = 14
$ws12.Range("C3").Value = 17.5

$ws12.Range("B7").Value = 18
$ws12.Range("C7").Value = 22.5

$ws12.Range("B8").Value = 19
$ws12.Range("C8").Value = 23.8

# --- Sheet: "Species qualification" ---
$ws4 = $wb.Worksheets.Item("Species qualification")

$ws4.Range("C3").Value = 228
$ws4.Range("C4").Value = 288
